$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp banner in A1
$ws.Range("A1").Value = "Datos actualizados a 19 de Junio de 2020 a las 18:18"

# Refresh country case numbers (new data pull) and re-rank rows whose
# position changed relative to neighbouring countries after the refresh.

# Row 4: Estados Unidos -> Estados Unidos
$ws.Cells.Item(4,2).Value = 2275701
$ws.Cells.Item(4,3).Value = 12050
$ws.Cells.Item(4,4).Value = 931499
$ws.Cells.Item(4,5).Value = 1223325
$ws.Cells.Item(4,7).Value = 189
$ws.Cells.Item(4,8).Value = 120877

# Row 7: India -> India
$ws.Cells.Item(7,2).Value = 385276
$ws.Cells.Item(7,3).Value = 4185
$ws.Cells.Item(7,4).Value = 208169
$ws.Cells.Item(7,5).Value = 164430
$ws.Cells.Item(7,7).Value = 73
$ws.Cells.Item(7,8).Value = 12677

# Row 12: Chile -> Chile
$ws.Cells.Item(12,2).Value = 231393
$ws.Cells.Item(12,3).Value = 6290
$ws.Cells.Item(12,5).Value = 40859
$ws.Cells.Item(12,7).Value = 252
$ws.Cells.Item(12,8).Value = 4093

# Row 36: Portugal -> Portugal
$ws.Cells.Item(36,2).Value = 38464
$ws.Cells.Item(36,3).Value = 375
$ws.Cells.Item(36,4).Value = 24477
$ws.Cells.Item(36,5).Value = 12460
$ws.Cells.Item(36,7).Value = 3
$ws.Cells.Item(36,8).Value = 1527

# Row 56: Kazajistan -> Kazajistan
$ws.Cells.Item(56,5).Value = 6043
$ws.Cells.Item(56,7).Value = 13
$ws.Cells.Item(56,8).Value = 113

# Row 57: Moldavia -> Ghana
$ws.Cells.Item(57,1).Value = "Ghana"
$ws.Cells.Item(57,2).Value = 13203
$ws.Cells.Item(57,3).Value = 274
$ws.Cells.Item(57,4).Value = 4548
$ws.Cells.Item(57,5).Value = 8585
$ws.Cells.Item(57,7).Value = 4
$ws.Cells.Item(57,8).Value = 70

# Row 58: Ghana -> Moldavia
$ws.Cells.Item(58,1).Value = "Moldavia"
$ws.Cells.Item(58,2).Value = 13106
$ws.Cells.Item(58,4).Value = 7525
$ws.Cells.Item(58,5).Value = 5132
$ws.Cells.Item(58,7).Value = 5
$ws.Cells.Item(58,8).Value = 449

# Row 63: Argelia -> Azerbaiyan
$ws.Cells.Item(63,1).Value = "Azerbaiyan"
$ws.Cells.Item(63,2).Value = 11767
$ws.Cells.Item(63,3).Value = 438
$ws.Cells.Item(63,4).Value = 6325
$ws.Cells.Item(63,5).Value = 5299
$ws.Cells.Item(63,7).Value = 4
$ws.Cells.Item(63,8).Value = 143

# Row 64: Azerbaiyan -> Argelia
$ws.Cells.Item(64,1).Value = "Argelia"
$ws.Cells.Item(64,2).Value = 11504
$ws.Cells.Item(64,3).Value = 119
$ws.Cells.Item(64,4).Value = 8196
$ws.Cells.Item(64,5).Value = 2483
$ws.Cells.Item(64,7).Value = 14
$ws.Cells.Item(64,8).Value = 825

# Row 67: Chequia -> Chequia
$ws.Cells.Item(67,2).Value = 10361
$ws.Cells.Item(67,3).Value = 81
$ws.Cells.Item(67,4).Value = 7472
$ws.Cells.Item(67,5).Value = 2554
$ws.Cells.Item(67,7).Value = 1
$ws.Cells.Item(67,8).Value = 335

# Row 87: Luxemburgo -> Luxemburgo
$ws.Cells.Item(87,2).Value = 4099
$ws.Cells.Item(87,3).Value = 8
$ws.Cells.Item(87,4).Value = 3944
$ws.Cells.Item(87,5).Value = 45

# Row 93: Grecia -> Grecia
$ws.Cells.Item(93,2).Value = 3237
$ws.Cells.Item(93,3).Value = 10
$ws.Cells.Item(93,5).Value = 1674
$ws.Cells.Item(93,7).Value = 1
$ws.Cells.Item(93,8).Value = 189

# Row 127: Jordania -> Jordania
$ws.Cells.Item(127,2).Value = 1008
$ws.Cells.Item(127,3).Value = 7
$ws.Cells.Item(127,4).Value = 708
$ws.Cells.Item(127,5).Value = 291

# Row 131: Georgia -> Georgia
$ws.Cells.Item(131,2).Value = 896
$ws.Cells.Item(131,3).Value = 3
$ws.Cells.Item(131,5).Value = 141

# Row 133: Principado de Andorra -> Republica del Chad
$ws.Cells.Item(133,1).Value = "Republica del Chad"
$ws.Cells.Item(133,2).Value = 858
$ws.Cells.Item(133,3).Value = 4
$ws.Cells.Item(133,4).Value = 742
$ws.Cells.Item(133,5).Value = 42
$ws.Cells.Item(133,8).Value = 74

# Row 134: Republica del Chad -> Principado de Andorra
$ws.Cells.Item(134,1).Value = "Principado de Andorra"
$ws.Cells.Item(134,2).Value = 855
$ws.Cells.Item(134,4).Value = 792
$ws.Cells.Item(134,5).Value = 11
$ws.Cells.Item(134,8).Value = 52

# Row 141: Estado de Palestina -> Mozambique
$ws.Cells.Item(141,1).Value = "Mozambique"
$ws.Cells.Item(141,2).Value = 668
$ws.Cells.Item(141,3).Value = 6
$ws.Cells.Item(141,4).Value = 177
$ws.Cells.Item(141,5).Value = 487
$ws.Cells.Item(141,8).Value = 4

# Row 142: Malta -> Estado de Palestina
$ws.Cells.Item(142,1).Value = "Estado de Palestina"
$ws.Cells.Item(142,3).Value = 63
$ws.Cells.Item(142,4).Value = 437
$ws.Cells.Item(142,5).Value = 223
$ws.Cells.Item(142,8).Value = 3

# Row 143: Mozambique -> Malta
$ws.Cells.Item(143,1).Value = "Malta"
$ws.Cells.Item(143,2).Value = 663
$ws.Cells.Item(143,4).Value = 613
$ws.Cells.Item(143,5).Value = 41
$ws.Cells.Item(143,8).Value = 9

# Row 202: Dominica -> Fiyi
$ws.Cells.Item(202,1).Value = "Fiyi"

# Row 203: Fiyi -> Dominica
$ws.Cells.Item(203,1).Value = "Dominica"

# Row 206: Islas Malvinas -> Groenlandia
$ws.Cells.Item(206,1).Value = "Groenlandia"

# Row 207: Groenlandia -> Islas Malvinas
$ws.Cells.Item(207,1).Value = "Islas Malvinas"

# Row 210: Seychelles -> Montserrat
$ws.Cells.Item(210,1).Value = "Montserrat"
$ws.Cells.Item(210,4).Value = 10
$ws.Cells.Item(210,8).Value = 1

# Row 211: Montserrat -> Seychelles
$ws.Cells.Item(211,1).Value = "Seychelles"
$ws.Cells.Item(211,4).Value = 11
$ws.Cells.Item(211,8).Value = 0

# Row 213: Papua Nueva Guinea -> Islas Virgenes Britanicas
$ws.Cells.Item(213,1).Value = "Islas Virgenes Britanicas"
$ws.Cells.Item(213,4).Value = 7
$ws.Cells.Item(213,8).Value = 1

# Row 214: Islas Virgenes Britanicas -> Papua Nueva Guinea
$ws.Cells.Item(214,1).Value = "Papua Nueva Guinea"
$ws.Cells.Item(214,4).Value = 8
$ws.Cells.Item(214,8).Value = 0
